# Scheduled runner update: refresh cached market-price / leve-profit figures
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
#  LeveProfitNQ/HQ) across each job sheet, per latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 399
$ws.Range("I2").Value = 288.6
$ws.Range("K2").Value = 288.6
$ws.Range("M2").Value = -175.6
$ws.Range("H18").Value = 866.6667
$ws.Range("I18").Value = 866.6667
$ws.Range("K18").Value = 866.6667
$ws.Range("M18").Value = -582.6667
$ws.Range("H74").Value = 6657.0303
$ws.Range("I74").Value = 5853.5557
$ws.Range("K74").Value = 5853.5557
$ws.Range("M74").Value = -4917.5557
$ws.Range("H77").Value = 6657.0303
$ws.Range("I77").Value = 5853.5557
$ws.Range("K77").Value = 29267.7785
$ws.Range("M77").Value = -24587.7785
$ws.Range("H86").Value = 4437.5
$ws.Range("I86").Value = 4360
$ws.Range("K86").Value = 4360
$ws.Range("M86").Value = -3237
$ws.Range("H89").Value = 4437.5
$ws.Range("I89").Value = 4360
$ws.Range("K89").Value = 21800
$ws.Range("M89").Value = -16184
$ws.Range("H92").Value = 1134.5
$ws.Range("I92").Value = 602.3077
$ws.Range("K92").Value = 602.3077
$ws.Range("M92").Value = 645.6923
$ws.Range("H96").Value = 438.9091
$ws.Range("I96").Value = 472.7
$ws.Range("K96").Value = 1418.1
$ws.Range("M96").Value = -45.09999999999991
$ws.Range("H98").Value = 687.2857
$ws.Range("I98").Value = 696.725
$ws.Range("K98").Value = 696.725
$ws.Range("M98").Value = 801.275
$ws.Range("H122").Value = 687.2857
$ws.Range("I122").Value = 696.725
$ws.Range("K122").Value = 2090.175
$ws.Range("M122").Value = 359.8249999999998
$ws.Range("H127").Value = 1196
$ws.Range("I127").Value = 1196
$ws.Range("K127").Value = 3588
$ws.Range("M127").Value = 1372
$ws.Range("H132").Value = 7032.25
$ws.Range("I132").Value = 1931.7333
$ws.Range("J132").Value = 22333.8
$ws.Range("K132").Value = 5795.199900000001
$ws.Range("L132").Value = 67001.39999999999
$ws.Range("M132").Value = -3265.199900000001
$ws.Range("N132").Value = -72061.39999999999
$ws.Range("H137").Value = 1581.4736
$ws.Range("I137").Value = 1318.625
$ws.Range("K137").Value = 3955.875
$ws.Range("M137").Value = -1405.875
$ws.Range("H138").Value = 2306.149
$ws.Range("J138").Value = 2438.862
$ws.Range("L138").Value = 7316.586
$ws.Range("N138").Value = -17596.586
$ws.Range("H141").Value = 2028.8096
$ws.Range("I141").Value = 2100.7
$ws.Range("J141").Value = 591
$ws.Range("K141").Value = 6302.099999999999
$ws.Range("L141").Value = 1773
$ws.Range("M141").Value = -1122.099999999999
$ws.Range("N141").Value = -12133

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8195.764999999999
$ws.Range("I61").Value = 7673.2085
$ws.Range("J61").Value = 9449.9
$ws.Range("K61").Value = 7673.2085
$ws.Range("L61").Value = 9449.9
$ws.Range("M61").Value = -7461.2085
$ws.Range("N61").Value = -9873.9
$ws.Range("H102").Value = 6750
$ws.Range("I102").Value = 4125
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 4125
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -2503
$ws.Range("N102").Value = -15244
$ws.Range("H122").Value = 4488.5835
$ws.Range("I122").Value = 4485.4287
$ws.Range("K122").Value = 13456.2861
$ws.Range("M122").Value = -11006.2861
$ws.Range("H136").Value = 8195.764999999999
$ws.Range("I136").Value = 7673.2085
$ws.Range("J136").Value = 9449.9
$ws.Range("K136").Value = 23019.6255
$ws.Range("L136").Value = 28349.7
$ws.Range("M136").Value = -20469.6255
$ws.Range("N136").Value = -33449.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1422.6666
$ws.Range("I105").Value = 1547.7142
$ws.Range("J105").Value = 985
$ws.Range("K105").Value = 1547.7142
$ws.Range("L105").Value = 985
$ws.Range("M105").Value = 199.2858000000001
$ws.Range("N105").Value = -4479

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2775.8125
$ws.Range("J16").Value = 4834.364
$ws.Range("L16").Value = 4834.364
$ws.Range("N16").Value = -5408.364
$ws.Range("H31").Value = 5295
$ws.Range("J31").Value = 5618.8184
$ws.Range("L31").Value = 5618.8184
$ws.Range("N31").Value = -6208.8184
$ws.Range("H34").Value = 5295
$ws.Range("J34").Value = 5618.8184
$ws.Range("L34").Value = 5618.8184
$ws.Range("N34").Value = -6022.8184
$ws.Range("H50").Value = 27500
$ws.Range("J50").Value = 27500
$ws.Range("L50").Value = 27500
$ws.Range("N50").Value = -28750
$ws.Range("H58").Value = 6146.9653
$ws.Range("I58").Value = 4069.6191
$ws.Range("J58").Value = 11600
$ws.Range("K58").Value = 4069.6191
$ws.Range("L58").Value = 11600
$ws.Range("M58").Value = -3866.6191
$ws.Range("N58").Value = -12006
$ws.Range("H62").Value = 7444.75
$ws.Range("I62").Value = 7333.3335
$ws.Range("J62").Value = 7779
$ws.Range("K62").Value = 7333.3335
$ws.Range("L62").Value = 7779
$ws.Range("M62").Value = -6709.3335
$ws.Range("N62").Value = -9027
$ws.Range("H65").Value = 7444.75
$ws.Range("I65").Value = 7333.3335
$ws.Range("J65").Value = 7779
$ws.Range("K65").Value = 36666.6675
$ws.Range("L65").Value = 38895
$ws.Range("M65").Value = -33546.6675
$ws.Range("N65").Value = -45135
$ws.Range("H113").Value = 2775.8125
$ws.Range("J113").Value = 4834.364
$ws.Range("L113").Value = 4834.364
$ws.Range("N113").Value = -9174.364
$ws.Range("H132").Value = 3375.3125
$ws.Range("I132").Value = 3020.4546
$ws.Range("J132").Value = 4156
$ws.Range("K132").Value = 9061.363799999999
$ws.Range("L132").Value = 12468
$ws.Range("M132").Value = -6531.363799999999
$ws.Range("N132").Value = -17528
$ws.Range("H134").Value = 5847.476
$ws.Range("I134").Value = 4599.8335
$ws.Range("J134").Value = 13333.333
$ws.Range("K134").Value = 13799.5005
$ws.Range("L134").Value = 39999.999
$ws.Range("M134").Value = -11264.5005
$ws.Range("N134").Value = -45069.999
$ws.Range("H136").Value = 6146.9653
$ws.Range("I136").Value = 4069.6191
$ws.Range("J136").Value = 11600
$ws.Range("K136").Value = 12208.8573
$ws.Range("L136").Value = 34800
$ws.Range("M136").Value = -9658.8573
$ws.Range("N136").Value = -39900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19232680
$ws.Range("I131").Value = 41667676
$ws.Range("J131").Value = 2682.9285
$ws.Range("K131").Value = 125003028
$ws.Range("L131").Value = 8048.7855
$ws.Range("M131").Value = -124997988
$ws.Range("N131").Value = -18128.7855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14012.15
$ws.Range("I70").Value = 11989.223
$ws.Range("K70").Value = 11989.223
$ws.Range("M70").Value = -11719.223
$ws.Range("H73").Value = 14012.15
$ws.Range("I73").Value = 11989.223
$ws.Range("K73").Value = 11989.223
$ws.Range("M73").Value = -11053.223
$ws.Range("H97").Value = 3754.2222
$ws.Range("I97").Value = 988.0909
$ws.Range("J97").Value = 8101
$ws.Range("K97").Value = 988.0909
$ws.Range("L97").Value = 8101
$ws.Range("M97").Value = -492.0909
$ws.Range("N97").Value = -9093
$ws.Range("H132").Value = 2006.3529
$ws.Range("I132").Value = 1350.5
$ws.Range("K132").Value = 4051.5
$ws.Range("M132").Value = -1521.5
$ws.Range("H137").Value = 61745
$ws.Range("I137").Value = 61745
$ws.Range("K137").Value = 61745
$ws.Range("M137").Value = -56645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1457.75
$ws.Range("I16").Value = 1365.15
$ws.Range("K16").Value = 1365.15
$ws.Range("M16").Value = -1195.15
$ws.Range("H22").Value = 1202.4166
$ws.Range("I22").Value = 892.9
$ws.Range("K22").Value = 892.9
$ws.Range("M22").Value = -597.9
$ws.Range("H27").Value = 1202.4166
$ws.Range("I27").Value = 892.9
$ws.Range("K27").Value = 892.9
$ws.Range("M27").Value = -785.9
$ws.Range("H46").Value = 13265.676
$ws.Range("I46").Value = 8335.117
$ws.Range("J46").Value = 17456.65
$ws.Range("K46").Value = 8335.117
$ws.Range("L46").Value = 17456.65
$ws.Range("M46").Value = -8147.117
$ws.Range("N46").Value = -17832.65
$ws.Range("H132").Value = 13818.565
$ws.Range("I132").Value = 14466.65
$ws.Range("J132").Value = 9498
$ws.Range("K132").Value = 43399.95
$ws.Range("L132").Value = 28494
$ws.Range("M132").Value = -40869.95
$ws.Range("N132").Value = -33554

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5798.1665
$ws.Range("I62").Value = 5092.467
$ws.Range("K62").Value = 5092.467
$ws.Range("M62").Value = -4468.467
$ws.Range("H65").Value = 5798.1665
$ws.Range("I65").Value = 5092.467
$ws.Range("K65").Value = 25462.335
$ws.Range("M65").Value = -22342.335
$ws.Range("H136").Value = 7662.7144
$ws.Range("I136").Value = 5422.625
$ws.Range("J136").Value = 10649.5
$ws.Range("K136").Value = 16267.875
$ws.Range("L136").Value = 31948.5
$ws.Range("M136").Value = -13717.875
$ws.Range("N136").Value = -37048.5
